$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 45, pushing existing rows 45-92 down to 47-94.
$ws.Rows.Item(45).Resize(2).Insert()

# Populate the two newly inserted rows (45 and 46) with new record data.
# Columns A,B,C,E,F,G,H,I,J,Q,T are constant for this data block.

# Row 45: Mandarina / Murcott / Segunda
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C45").Value = "Arica y Parinacota"
$ws.Range("D45").Value = (Get-Date -Year 2022 -Month 3 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100102
$ws.Range("H45").Value = "Cítricos"
$ws.Range("I45").Value = 100102004
$ws.Range("J45").Value = "Mandarina"
$ws.Range("K45").Value = "Murcott"
$ws.Range("L45").Value = "Segunda"
$ws.Range("M45").Value = 125
$ws.Range("N45").Value = 13000
$ws.Range("O45").Value = 13000
$ws.Range("P45").Value = 13000
$ws.Range("Q45").Value = "$/caja 20 kilos"
$ws.Range("R45").Value = "Región de Coquimbo"
$ws.Range("S45").Value = 650
$ws.Range("T45").Value = 20

# Row 46: Mandarina / Murcott / Tercera
$ws.Range("A46").Value = 1
$ws.Range("B46").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C46").Value = "Arica y Parinacota"
$ws.Range("D46").Value = (Get-Date -Year 2022 -Month 3 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E46").Value = 15
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100102
$ws.Range("H46").Value = "Cítricos"
$ws.Range("I46").Value = 100102004
$ws.Range("J46").Value = "Mandarina"
$ws.Range("K46").Value = "Murcott"
$ws.Range("L46").Value = "Tercera"
$ws.Range("M46").Value = 125
$ws.Range("N46").Value = 14000
$ws.Range("O46").Value = 14000
$ws.Range("P46").Value = 14000
$ws.Range("Q46").Value = "$/caja 20 kilos"
$ws.Range("R46").Value = "Región de Coquimbo"
$ws.Range("S46").Value = 700
$ws.Range("T46").Value = 20
